# Added Contact Module test as well
# Remove the "title" column (with Mr./Dr./Mrs. values) from the "contacts"
# sheet, shifting firstname/lastname/company left by one column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")

# Delete entire column A (title / Mr. / Dr. / Mrs.) and shift remaining
# columns (firstname, lastname, company) to the left.
$ws.Columns.Item(1).Delete()

# Select D2 to mirror the active-cell selection recorded after the edit.
$ws.Activate()
$ws.Range("D2").Select()

$wb.Save()
